$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new hourly ranking rows were inserted right before the existing
# "2026/12/29" block (previously starting at row 711), shifting all rows
# from the old row 711 through 752 down by two (now rows 713-754).
$ws.Rows.Item(711).Insert()
$ws.Rows.Item(711).Insert()

# Populate the two newly inserted rows with their data. The date column
# is forced to plain text (via a temporary Text number format, cleared
# immediately after entry) so Excel doesn't auto-convert the "yyyy/mm/dd"
# strings into date serial numbers, keeping it consistent with the other
# inline-string date cells in the sheet.
$ws.Range("A711").NumberFormat = "@"
$ws.Range("A711").Value = "2026/01/23"
$ws.Range("A711").ClearFormats()
$ws.Range("B711").Value = "金"
$ws.Range("C711").Value = 23
$ws.Range("D711").Value = 161

$ws.Range("A712").NumberFormat = "@"
$ws.Range("A712").Value = "2026/01/24"
$ws.Range("A712").ClearFormats()
$ws.Range("B712").Value = "土"
$ws.Range("C712").Value = 2
$ws.Range("D712").Value = 166
